# Add season record columns (Wins, Losses, Ties) to the player table.
# The previous export only captured team statistics, not the season record,
# so every row of the sheet gets the Braves' 1998 record appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, styled the same as the rest of row 1 (bold, centered,
# thin border all around - matches style index "1" used by the existing
# header cells such as AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Season record (106 wins, 56 losses, 0 ties) repeated for every player row.
$wins = 106
$losses = 56
$ties = 0

$lastRow = 44
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins
    $ws.Cells.Item($r, 31).Value = $losses
    $ws.Cells.Item($r, 32).Value = $ties
}
